$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I (I0) and J (IF), copying the H1 header format
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-20
$data = @{
    2  = @(1, 3)
    3  = @(5, 6)
    4  = @(1, 4)
    5  = @(8, 8)
    6  = @(8, 9)
    7  = @(2, 4)
    8  = @(3, 4)
    9  = @(5, 7)
    10 = @(5, 6)
    11 = @(1, 3)
    12 = @(1, 2)
    13 = @(1, 4)
    14 = @(8, 8)
    15 = @(6, 7)
    16 = @(8, 8)
    17 = @(6, 6)
    18 = @(7, 7)
    19 = @(5, 5)
    20 = @(7, 8)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
